$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price column (C) for rows 3-34 ---
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '$129'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '$289'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '$224'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '$110'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = '$199'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '$360'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '$124'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = '$339'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '$179'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '$202'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '$268'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '$214'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '$229'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '$78'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = '$369'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '$172'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '$199'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '$360'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = '$469'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '$119'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '$229'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = '$315'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = '$76'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = '$144'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = '$154'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = '$299'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = '$405'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = '$201'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = '$3,129'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = '$164'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = '$549'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = '$1,459'

# --- Rows 35-38 shift up (CPU model/specs) + new row 38 ---
# Row 35: Core i3-10105
$ws.Range("B35").Value = 'Core i3-10105'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = '$59'
$ws.Range("D35").Value = 'LGA 1200'
$ws.Range("E35").Value = 'Quad-Core'
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = '8'
$ws.Range("G35").Value = '3.7 GHz'
$ws.Range("H35").Value = '4.4 GHz'

# Row 36: Core i9-10900K
$ws.Range("B36").Value = 'Core i9-10900K'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = '$179'
$ws.Range("D36").Value = 'LGA 1200'
$ws.Range("E36").Value = '10-Core'
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = '20'
$ws.Range("G36").Value = '3.7 GHz'
$ws.Range("H36").Value = '5.30 GHz'

# Row 37: Core i9-12900
$ws.Range("B37").Value = 'Core i9-12900'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = '$359'
$ws.Range("D37").Value = 'LGA 1700'
$ws.Range("E37").Value = '16-Core (8P+8E)'
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = '24'
$ws.Range("G37").Value = 'P-core Base Frequency: 2.4 GHzE-core Base Frequency: 1.8 GHz'
$ws.Range("H37").Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 5.1 GHzP-core Turbo Frequency: Up to 5.0 GHzE-core Turbo Frequency: Up to 3.8 GHz'

# Row 38: Core i7-12700
$ws.Range("B38").Value = 'Core i7-12700'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = '$354'
$ws.Range("D38").Value = 'LGA 1700'
$ws.Range("E38").Value = '12-Core (8P+4E)'
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = '20'
$ws.Range("G38").Value = 'P-core Base Frequency: 2.1 GHzE-core Base Frequency: 1.6 GHz'
$ws.Range("H38").Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 4.9 GHzP-core Turbo Frequency: Up to 4.8 GHzE-core Turbo Frequency: Up to 3.6 GHz'

